$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign as literal text. A leading apostrophe forces Excel to
# keep numeric-looking strings like "8.00" / "1.00" / "0.117" as text
# instead of silently converting them to numbers (which would drop
# trailing zeros / thousands-style separators and change the cell type).
$apost = "'"
function Set-TextValue($cellRef, [string]$value) {
    $ws.Range($cellRef).Value = $apost + $value
}

Set-TextValue "D2" '63.853.81'
Set-TextValue "E2" '  +0.96%  '

Set-TextValue "D3" '3.323.96'
Set-TextValue "E3" '  +2.01%  '

Set-TextValue "E4" '  +0.07%  '

Set-TextValue "E5" '  +1.64%  '

Set-TextValue "D6" '142.89'
Set-TextValue "E6" '  +1.32%  '

Set-TextValue "E7" '  +0.03%  '

Set-TextValue "D8" '3.323.72'
Set-TextValue "E8" '  +2.12%  '

Set-TextValue "E9" '  -0.15%  '

Set-TextValue "E10" '  +1.54%  '

Set-TextValue "D11" '5.54'
Set-TextValue "E11" '  +4.01%  '

Set-TextValue "D12" '0.469'
Set-TextValue "E12" '  +0.74%  '

Set-TextValue "E13" '  +0.14%  '

Set-TextValue "D14" '35.09'
Set-TextValue "E14" '  +1.55%  '

Set-TextValue "D15" '3.875.56'
Set-TextValue "E15" '  +2.38%  '

Set-TextValue "E16" '  +0.35%  '

Set-TextValue "D17" '3.325.17'
Set-TextValue "E17" '  +2.44%  '

Set-TextValue "D18" '63.930.61'
Set-TextValue "E18" '  +0.96%  '

Set-TextValue "E19" '  +1.40%  '

Set-TextValue "D20" '480.49'
Set-TextValue "E20" '  +0.91%  '

Set-TextValue "D21" '14.13'
Set-TextValue "E21" '  -0.24%  '

Set-TextValue "E22" '  +1.64%  '

Set-TextValue "D23" '8.00'
Set-TextValue "E23" '  +0.60%  '

Set-TextValue "E24" '  +4.04%  '

Set-TextValue "D25" '84.89'
Set-TextValue "E25" '  +1.36%  '

Set-TextValue "E27" '  +1.66%  '

Set-TextValue "B28" 'FirstDigitalUSD'
Set-TextValue "C28" 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D28" '1.00'
Set-TextValue "E28" '  +0.03%  '

Set-TextValue "B29" 'RenderToken'
Set-TextValue "C29" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D29" '8.27'
Set-TextValue "E29" '  +2.00%  '

Set-TextValue "D30" '7.18'
Set-TextValue "E30" '  -2.73%  '

Set-TextValue "E31" '  +1.26%  '

Set-TextValue "D32" '28.96'
Set-TextValue "E32" '  +5.05%  '

Set-TextValue "E33" '  -1.16%  '

Set-TextValue "D34" '2.53'
Set-TextValue "E34" '  -0.66%  '

Set-TextValue "E35" '  +0.25%  '

Set-TextValue "D36" '6.07'
Set-TextValue "E36" '  +2.98%  '

Set-TextValue "D37" '0.0₃0752'
Set-TextValue "E37" '  +4.17%  '

Set-TextValue "D38" '52.42'
Set-TextValue "E38" '  -0.92%  '

Set-TextValue "D39" '0.0399'
Set-TextValue "E39" '  +1.47%  '

Set-TextValue "D40" '3.123.31'
Set-TextValue "E40" '  +4.00%  '

Set-TextValue "D41" '431.18'
Set-TextValue "E41" '  +2.60%  '

Set-TextValue "E42" '  -0.50%  '

Set-TextValue "B43" 'Kaspa'
Set-TextValue "C43" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D43" '0.117'
Set-TextValue "E43" '  +5.79%  '

Set-TextValue "B44" 'Cosmos'
Set-TextValue "C44" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D44" '8.35'
Set-TextValue "E44" '  -0.60%  '

Set-TextValue "D45" '0.266'
Set-TextValue "E45" '  -0.91%  '

Set-TextValue "D46" '2.25'
Set-TextValue "E46" '  +3.75%  '

Set-TextValue "B47" 'InjectiveProtocol'
Set-TextValue "C47" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D47" '26.43'
Set-TextValue "E47" '  +1.90%  '

Set-TextValue "B48" 'Arweave'
Set-TextValue "C48" 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue "D48" '36.17'
Set-TextValue "E48" '  +8.31%  '

Set-TextValue "E50" '  -0.37%  '

Set-TextValue "E51" '  -0.65%  '
